# Apply the "Saldo" roster update described by the commit diff.
#
# Changes (top of sheet "Export"):
#   - remove account 005642649 / VR / 500000
#   - remove account 004526450 / MSD / 32000
#   - insert account 004460487 / PEDRO   / 15955.08 (just above CATARINE)
#   - insert account 004550605 / REJANE  / 10473.55 (just above ADELE)
# Changes (far down the sheet, near balances around 0.8):
#   - remove the now-duplicate account 004460487 / PEDRO / 0.83 row
#     (PEDRO's current balance moved to the top of the sheet instead)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Find-RowByAccount($sheet, [string]$account, [double]$balance) {
    $used = $sheet.UsedRange
    $lastRow = $used.Rows.Count
    for ($r = 1; $r -le $lastRow; $r++) {
        $a = $sheet.Cells.Item($r, 1).Value2
        if ($a -eq $account) {
            $c = $sheet.Cells.Item($r, 3).Value2
            if ([math]::Abs([double]$c - $balance) -lt 0.005) {
                return $r
            }
        }
    }
    return -1
}

# 1) Delete the stale low-balance PEDRO row further down the sheet first
#    (processing the bottom-most change first keeps the row numbers for
#    the edits above it stable).
$rowPedroOld = Find-RowByAccount $ws "004460487" 0.83
if ($rowPedroOld -gt 0) {
    $ws.Rows.Item($rowPedroOld).Delete()
}

# 2) Insert REJANE just above ADELE (004575632)
$rowAdele = Find-RowByAccount $ws "004575632" 7057.33
$ws.Rows.Item($rowAdele).Insert()
$ws.Cells.Item($rowAdele, 1).NumberFormat = "@"
$ws.Cells.Item($rowAdele, 1).Value = "004550605"
$ws.Cells.Item($rowAdele, 2).Value = "REJANE"
$ws.Cells.Item($rowAdele, 3).Value = 10473.55

# 3) Insert PEDRO (new balance) just above CATARINE (004693349)
$rowCatarine = Find-RowByAccount $ws "004693349" 15398.41
$ws.Rows.Item($rowCatarine).Insert()
$ws.Cells.Item($rowCatarine, 1).NumberFormat = "@"
$ws.Cells.Item($rowCatarine, 1).Value = "004460487"
$ws.Cells.Item($rowCatarine, 2).Value = "PEDRO"
$ws.Cells.Item($rowCatarine, 3).Value = 15955.08

# 4) Remove the VR / 500000 row
$rowVr = Find-RowByAccount $ws "005642649" 500000
if ($rowVr -gt 0) {
    $ws.Rows.Item($rowVr).Delete()
}

# 5) Remove the MSD / 32000 row
$rowMsd = Find-RowByAccount $ws "004526450" 32000
if ($rowMsd -gt 0) {
    $ws.Rows.Item($rowMsd).Delete()
}
